$wb = $excel.ActiveWorkbook

# --- Sheet "Metadata" (sheet1) ---
$ws1 = $wb.Worksheets.Item("Metadata")

# Version 5.0.0 -> 6.0.0
$ws1.Range("B3").Value = "6.0.0"

# Date updated
$ws1.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher value was blank, now "Alvearie Team"
$ws1.Range("B9").Value = "Alvearie Team"

# Row 10 was "Contact" / "No display for ContactDetail" -> becomes "Jurisdiction" / "United States of America"
$ws1.Range("A10").Value = "Jurisdiction"
$ws1.Range("B10").Value = "United States of America"

# Row 11 was a duplicate "Contact" / "No display for ContactDetail" row -> delete it entirely
$ws1.Rows.Item(11).Delete()

# --- Sheet "Elements" (sheet2) ---
$ws2 = $wb.Worksheets.Item("Elements")

# Root Extension row: Short / Definition columns get the specific extension title & description
# instead of the generic "Extension" / "An Extension" placeholders.
$ws2.Range("K2").Value = "ACA Health Insurance Oversight System Product"
$ws2.Range("L2").Value = "Code for Affordable Care Act (ACA) HIOS product ID and component code of the associated plan"

# Column K ("Short") grew wider text, so its best-fit width increases.
$ws2.Columns.Item(11).ColumnWidth = 45.546875
